$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SANDRA's balance (row 2, column C) from 110344.98 to 95201.98
$ws.Cells.Item(2, 3).Value = 95201.98

# Delete the row for account 004574428 / GUILHERME (row 6), shifting subsequent rows up
$ws.Rows.Item(6).Delete()
